# Generate Report for Handoff
# Updates the handoff timestamps for the e05b360d-2fd3-4275-9404-14655c46e33e.md
# row on the Overview sheet and both language sheets (zh-cn, de-de) to reflect
# a newly generated handoff xliff.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-11-02 03:48:20"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-11-02 03:48:07"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-11-02 03:48:20"
